$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: update date, price min/max/avg and $/Kg price
$ws.Cells.Item(9, 4).Value = 44627
$ws.Cells.Item(9, 11).Value = 4000
$ws.Cells.Item(9, 12).Value = 4500
$ws.Cells.Item(9, 13).Value = 4250
$ws.Cells.Item(9, 16).Value = 71

# Row 10: update date and volume
$ws.Cells.Item(10, 4).Value = 44362
$ws.Cells.Item(10, 10).Value = 120

# New row 11: append a new record (same shape as the former row 10)
$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(11, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(11, 4).Value = 44421
$ws.Cells.Item(11, 4).NumberFormat = $ws.Cells.Item(10, 4).NumberFormat
$ws.Cells.Item(11, 5).Value = 15
$ws.Cells.Item(11, 6).Value = 100112001
$ws.Cells.Item(11, 7).Value = "Berenjena"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 100
$ws.Cells.Item(11, 11).Value = 8000
$ws.Cells.Item(11, 12).Value = 9000
$ws.Cells.Item(11, 13).Value = 8500
$ws.Cells.Item(11, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(11, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(11, 16).Value = 142
$ws.Cells.Item(11, 17).Value = 60
$ws.Cells.Item(11, 18).Value = "Hortaliza"
